$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (GatesS)
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 0.4657814258739751
$ws.Range("D3").Value = 0.1341093468370086
$ws.Range("E3").Value = "coldread_stopwatch_wpm"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1.47e-26"

# Row 4 (GatesT)
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 0.5804708609886242
$ws.Range("D4").Value = 0.1313706372987473
$ws.Range("E4").Value = "qa_coverage_line_%"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2.82e-09"
